$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.806.08'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.266.56'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '305.05'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '92.87'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('E7').Value = '  -0.88%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.485'
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '32.68'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.67'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.619.16'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.268.86'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('E17').Value = '  +3.18%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '41.747.17'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.92'
$ws.Range('E19').Value = '  +5.81%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0911'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.97'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('E22').Value = '  +0.99%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '243.53'
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.94'
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.97'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.62'
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('E29').Value = '  -5.30%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '34.64'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '159.18'
$ws.Range('E31').Value = '  +0.40%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.36'
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0742'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('E35').Value = '  -2.16%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.08'
$ws.Range('E36').Value = '  +2.74%  '
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('E40').Value = '  -1.08%  '
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.010.46'
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '19.53'
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.25'
$ws.Range('E44').Value = '  +12.75%  '
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.27'
$ws.Range('E46').Value = '  +1.63%  '
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '53.49'
$ws.Range('E48').Value = '  +2.97%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '73.13'
$ws.Range('E49').Value = '  +3.27%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.50'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.14'
$ws.Range('E51').Value = '  -0.23%  '
